$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / blank cell updates: (row, col, value)
$normalEdits = @(
    @(6,11,'x'),
    @(6,15,''),
    @(6,20,'StateArea'),
    @(6,21,'State Area'),
    @(6,22,'Area (in square meters) of state'),
    @(6,24,'Tiger/Line 2018'),
    @(6,25,'Tiger/Line 2018 Shapefiles'),
    @(7,11,''),
    @(7,15,'x'),
    @(7,20,'ParkArea'),
    @(7,21,'Park Area'),
    @(7,22,'Area (in square meters) of park or green space in a state).'),
    @(7,24,'OSM'),
    @(7,25,'OpenStreetMap'),
    @(54,20,'CntNaltT'),
    @(57,20,'AvNaltTime'),
    @(60,20,'PctNaltT'),
    @(67,20,'NaltCtTmDr'),
    @(68,20,'NaltCtTmBk'),
    @(69,20,'NaltCtTmWk'),
    @(76,20,'NaltAvTmDr'),
    @(77,20,'NaltAvTmBk'),
    @(78,20,'NaltAvTmWk'),
    @(85,20,'NaltTmDrP'),
    @(86,20,'NaltTmBkP'),
    @(87,20,'NaltTmWkP'),
    @(141,20,'AnyNalxDt'),
    @(142,20,'NalxPrStDt'),
    @(143,20,'NalxPresDt'),
    @(144,20,'AnyNalxFr'),
    @(145,20,'NalxPrStFr'),
    @(146,20,'NalxPresFr'),
    @(157,11,'x'),
    @(157,12,''),
    @(157,13,''),
    @(157,14,''),
    @(157,15,''),
    @(157,20,'CrrctExpS'),
    @(157,21,'Corrections Expenditures'),
    @(157,22,'Expenditures on corrections system and operation by the State alone'),
    @(157,24,'US Census, 2018'),
    @(157,25,'U.S. Census Bureau Annual Survey of State and Local Government Finances via Urban Institute & Tax Policy Center''s State and Local Finance Data Finder'),
    @(157,29,'Local data outside of police and fire expenditures was not available for Washington, D.C.'),
    @(158,11,'x'),
    @(158,12,''),
    @(158,13,''),
    @(158,14,''),
    @(158,15,''),
    @(158,20,'PlcFyrExpS'),
    @(158,21,'Police & Fire Expenditures'),
    @(158,22,'Expenditures on police and fire protection by the State alone'),
    @(158,24,'US Census, 2018'),
    @(158,25,'U.S. Census Bureau Annual Survey of State and Local Government Finances via Urban Institute & Tax Policy Center''s State and Local Finance Data Finder'),
    @(158,29,'Local data outside of police and fire expenditures was not available for Washington, D.C.'),
    @(159,11,'x'),
    @(159,12,''),
    @(159,13,''),
    @(159,14,''),
    @(159,15,''),
    @(159,20,'HlthExpS'),
    @(159,21,'Public Health Expenditures'),
    @(159,22,'Expenditures on public health and hospitals by the State alone'),
    @(159,24,'US Census, 2018'),
    @(159,25,'U.S. Census Bureau Annual Survey of State and Local Government Finances via Urban Institute & Tax Policy Center''s State and Local Finance Data Finder'),
    @(159,29,'Local data outside of police and fire expenditures was not available for Washington, D.C.'),
    @(160,11,'x'),
    @(160,12,''),
    @(160,13,''),
    @(160,14,''),
    @(160,15,''),
    @(160,20,'WlfrExpS'),
    @(160,21,'Public Welfare Expenditures'),
    @(160,22,'Expenditures on public welfare progrmas by the State alone'),
    @(160,24,'US Census, 2018'),
    @(160,25,'U.S. Census Bureau Annual Survey of State and Local Government Finances via Urban Institute & Tax Policy Center''s State and Local Finance Data Finder'),
    @(160,29,'Local data outside of police and fire expenditures was not available for Washington, D.C.'),
    @(161,20,'CrrctExpL'),
    @(161,21,'Corrections Expenditures (Local)'),
    @(161,22,'Expenditures on corrections system and operation by local governments alone'),
    @(162,20,'PlcFyrExpL'),
    @(162,21,'Police & Fire Expenditures (Local)'),
    @(162,22,'Expenditures on police and fire protection by the local government alone'),
    @(163,20,'HlthExpL'),
    @(163,21,'Public Health Expenditures (Local)'),
    @(163,22,'Expenditures on public health and hospitals by the local government alone'),
    @(164,20,'WlfrExpL'),
    @(164,21,'Public Welfare Expenditures (Local)'),
    @(164,22,'Expenditures on public welfare progrmas by the local government alone'),
    @(165,20,'CrrctExpT'),
    @(165,21,'Total Corrections Expenditures'),
    @(165,22,'Total expenditures on corrections system and operations'),
    @(166,20,'PlcFyrExpT'),
    @(166,21,'Total Police & Fire Expenditures'),
    @(166,22,'Total expenditures on police and fire protection'),
    @(167,20,'HlthExpT'),
    @(167,21,'Total Public Health & Welfare Expenditures'),
    @(167,22,'Total expenditures on public health and welfare'),
    @(168,20,'WlfrExpT'),
    @(168,21,'Total Public Welfare Expenditures'),
    @(168,22,'Total expenditures on public welfare programs'),
    @(169,11,''),
    @(169,12,'x'),
    @(169,13,'x'),
    @(169,14,'x'),
    @(169,15,'x'),
    @(169,20,'CrrctExp'),
    @(169,21,'Corrections expenditures'),
    @(169,24,'State and Local Finance Data Finder'),
    @(169,25,'State and Local Finance Data Finder'),
    @(169,29,''),
    @(170,11,''),
    @(170,12,'x'),
    @(170,13,'x'),
    @(170,14,'x'),
    @(170,15,'x'),
    @(170,20,'HlthExp'),
    @(170,21,'Public health expenditures'),
    @(170,22,'Total expenditures on public health and hospitals'),
    @(170,24,'State and Local Finance Data Finder'),
    @(170,25,'State and Local Finance Data Finder'),
    @(170,29,''),
    @(171,11,''),
    @(171,12,'x'),
    @(171,13,'x'),
    @(171,14,'x'),
    @(171,15,'x'),
    @(171,20,'PlcFyrExp'),
    @(171,21,'Police & fire expenditures'),
    @(171,22,'Total expenditures on police and fire protection'),
    @(171,24,'State and Local Finance Data Finder'),
    @(171,25,'State and Local Finance Data Finder'),
    @(171,29,''),
    @(172,11,''),
    @(172,12,'x'),
    @(172,13,'x'),
    @(172,14,'x'),
    @(172,15,'x'),
    @(172,20,'WlfrExp'),
    @(172,21,'Public welfare expenditures'),
    @(172,24,'State and Local Finance Data Finder'),
    @(172,25,'State and Local Finance Data Finder'),
    @(172,29,''),
    @(179,6,'x'),
    @(179,7,'x'),
    @(179,8,'x'),
    @(179,9,'x'),
    @(179,11,'x'),
    @(179,12,'x'),
    @(179,13,'x'),
    @(179,14,'x'),
    @(179,15,'x'),
    @(179,20,'HcvD'),
    @(179,21,'Hepatitis C Deaths'),
    @(179,22,'Total Hepatitis C deaths'),
    @(179,27,'integer'),
    @(180,6,'x'),
    @(180,7,'x'),
    @(180,8,'x'),
    @(180,9,'x'),
    @(180,11,'x'),
    @(180,12,'x'),
    @(180,13,'x'),
    @(180,14,'x'),
    @(180,15,'x'),
    @(180,20,'MlHcvD'),
    @(180,21,'Hepatitis C Deaths - Men'),
    @(180,22,'Hepatitis C deaths among men'),
    @(180,27,'integer'),
    @(181,6,'x'),
    @(181,7,'x'),
    @(181,8,'x'),
    @(181,9,'x'),
    @(181,11,'x'),
    @(181,12,'x'),
    @(181,13,'x'),
    @(181,14,'x'),
    @(181,15,'x'),
    @(181,20,'FlHcvD'),
    @(181,21,'Hepatitis C Deaths - Women'),
    @(181,22,'Hepatitis C deaths among women'),
    @(181,27,'integer'),
    @(182,6,'x'),
    @(182,7,'x'),
    @(182,8,'x'),
    @(182,9,'x'),
    @(182,11,'x'),
    @(182,12,'x'),
    @(182,13,'x'),
    @(182,14,'x'),
    @(182,15,'x'),
    @(182,20,'AmInHcvD'),
    @(182,21,'Hepatitis C Deaths - American Indian'),
    @(182,22,'Hepatitis C deaths among American Indian populations'),
    @(182,27,'integer'),
    @(183,6,'x'),
    @(183,7,'x'),
    @(183,8,'x'),
    @(183,9,'x'),
    @(183,20,'AsPiHcvD'),
    @(183,21,'Hepatitis C Deaths - Asian & Pacific Islander'),
    @(183,22,'Hepatitis C deaths among Asian and Pacific Islander populations'),
    @(183,27,'integer'),
    @(184,6,'x'),
    @(184,7,'x'),
    @(184,8,'x'),
    @(184,9,'x'),
    @(184,11,'x'),
    @(184,12,'x'),
    @(184,13,'x'),
    @(184,14,'x'),
    @(184,15,'x'),
    @(184,20,'BlkHcvD'),
    @(184,21,'Hepatitis C Deaths - Black'),
    @(184,22,'Hepatitis C deaths among Black population'),
    @(184,27,'integer'),
    @(185,6,'x'),
    @(185,7,'x'),
    @(185,8,'x'),
    @(185,9,'x'),
    @(185,11,'x'),
    @(185,12,'x'),
    @(185,13,'x'),
    @(185,14,'x'),
    @(185,15,'x'),
    @(185,20,'HspHcvD'),
    @(185,21,'Hepatitis C Deaths - Hispanic'),
    @(185,22,'Hepatitis C deaths among hispanic populations'),
    @(185,27,'integer'),
    @(186,6,'x'),
    @(186,7,'x'),
    @(186,8,'x'),
    @(186,9,'x'),
    @(186,11,'x'),
    @(186,12,'x'),
    @(186,13,'x'),
    @(186,14,'x'),
    @(186,15,'x'),
    @(186,20,'U50HcvD'),
    @(186,21,'Hepatitis C Deaths - Under 50 years old'),
    @(186,22,'Hepatitis C deaths in populations under 50 years of age'),
    @(186,27,'integer'),
    @(187,20,'A50_74HcvD'),
    @(187,21,'Hepatitis C Deaths - 50 to 74 years old'),
    @(187,22,'Hepatitis C deaths among populations between 50 and 74 years of age'),
    @(188,20,'O75HcvD'),
    @(188,21,'Hepatitis C Deaths - Over 75 years old'),
    @(188,22,'Hepatitis C deaths among populations over 75 years of age'),
    @(189,6,''),
    @(189,7,''),
    @(189,8,''),
    @(189,9,''),
    @(189,11,''),
    @(189,12,''),
    @(189,13,''),
    @(189,14,''),
    @(189,15,''),
    @(189,20,'TotHcv'),
    @(189,21,'Yearly Hepatitis C cases (2013-2016)'),
    @(189,22,'Mean total yearly Hepitatis C cases from 2013-2016'),
    @(189,27,'number'),
    @(190,6,''),
    @(190,7,''),
    @(190,8,''),
    @(190,9,''),
    @(190,11,''),
    @(190,12,''),
    @(190,13,''),
    @(190,14,''),
    @(190,15,''),
    @(190,20,'MlHcv'),
    @(190,21,'Yearly Hepatitis C cases - Men (2013-2016)'),
    @(190,22,'Mean yearly Hepatitis C cases in men from 2013-2016'),
    @(190,27,'number'),
    @(191,6,''),
    @(191,7,''),
    @(191,8,''),
    @(191,9,''),
    @(191,20,'FmHcv'),
    @(191,21,'Yearly Hepatitis C cases - Women (2013-2016)'),
    @(191,22,'Mean yearly Hepatitis C cases in women from 2013-2016'),
    @(191,27,'number'),
    @(192,6,''),
    @(192,7,''),
    @(192,8,''),
    @(192,9,''),
    @(192,11,''),
    @(192,12,''),
    @(192,13,''),
    @(192,14,''),
    @(192,15,''),
    @(192,20,'Un50Hcv'),
    @(192,21,'Yearly Hepatitis C cases - Under 50 years old (2013-2016)'),
    @(192,22,'Mean yearly Hepatatis C cases in people under 50 years of age from 2013-2016'),
    @(192,27,'number'),
    @(193,6,''),
    @(193,7,''),
    @(193,8,''),
    @(193,9,''),
    @(193,11,''),
    @(193,12,''),
    @(193,13,''),
    @(193,14,''),
    @(193,15,''),
    @(193,20,'A50_74Hcv'),
    @(193,21,'Yearly Hepatitis C cases - 50 to 74 years old (2013-2016)'),
    @(193,22,'Mean yearly Hepatitis C cases in people between 50 to 74 years of age from 2013-2016'),
    @(193,27,'number'),
    @(194,6,''),
    @(194,7,''),
    @(194,8,''),
    @(194,9,''),
    @(194,11,''),
    @(194,12,''),
    @(194,13,''),
    @(194,14,''),
    @(194,15,''),
    @(194,20,'Ov75Hcv'),
    @(194,21,'Yearly Hepatitis C cases - Over 75 years old (2013-2016)'),
    @(194,22,'Mean yearly Hepatitis C cases in people over 75 years of age from 2013-2016'),
    @(194,27,'number'),
    @(195,6,''),
    @(195,7,''),
    @(195,8,''),
    @(195,9,''),
    @(195,11,''),
    @(195,12,''),
    @(195,13,''),
    @(195,14,''),
    @(195,15,''),
    @(195,20,'BlkHcv'),
    @(195,21,'Yearly Hepatitis C cases - Black (2013-2016)'),
    @(195,22,'Mean yearly Hepatitis C cases in populations identified as non-hispanic Black alone across 2013-2016'),
    @(195,27,'number'),
    @(196,6,''),
    @(196,7,''),
    @(196,8,''),
    @(196,9,''),
    @(196,11,''),
    @(196,12,''),
    @(196,13,''),
    @(196,14,''),
    @(196,15,''),
    @(196,20,'NonBlkHcv'),
    @(196,21,'Yearly Hepatitis C cases - non-Black (2013-2016)'),
    @(196,22,'Mean yearly Hepatitis C cases in populations non-Black other race/ethnicity populations 2013-2016'),
    @(196,27,'number'),
    @(206,11,'x'),
    @(206,12,'x'),
    @(206,13,'x'),
    @(206,14,'x'),
    @(206,20,'AsHcvD'),
    @(206,21,'Hepatitis C deaths among Asian populations '),
    @(206,22,''),
    @(207,20,'NhPiHcvD'),
    @(207,21,'Hepatitis C deaths among Native Hawaiian and Pacific Islander populations '),
    @(208,20,'WhtHcvD'),
    @(208,21,'Hepatitis C deaths among White populations'),
    @(209,20,'MulHcvD'),
    @(209,21,'Hepatitis C deaths among Multiple Race populations '),
    @(210,11,''),
    @(210,12,''),
    @(210,13,''),
    @(210,14,''),
    @(210,20,'AvHcvD'),
    @(210,21,'Average Hepitatis C virus Deaths'),
    @(210,22,'Mean total yearly Hepatitis C deaths from 2013-2017, 2018-2022')
)

foreach ($edit in $normalEdits) {
    $ws.Cells.Item($edit[0], $edit[1]).Value = $edit[2]
}

# Numeric-looking text values: force text storage so they are written as
# strings (matching the source workbook's inline-string "Example" column)
# rather than being auto-converted to numbers by Excel.
$numericTextEdits = @(
    @(6,28,'4486028684.18'),
    @(7,28,'39296898972.47'),
    @(157,28,'915105'),
    @(158,28,'305035'),
    @(159,28,'915105'),
    @(160,28,'1830210'),
    @(161,28,'2196606'),
    @(162,28,'3385247'),
    @(163,28,'2196606'),
    @(164,28,'4393212'),
    @(165,28,'3111711'),
    @(166,28,'3690282'),
    @(167,28,'3111711'),
    @(168,28,'6223422'),
    @(169,28,'1897193'),
    @(170,28,'1897193'),
    @(171,28,'1897193'),
    @(172,28,'1897193'),
    @(179,28,'592'),
    @(180,28,'417'),
    @(181,28,'175'),
    @(182,28,'20'),
    @(183,28,'29'),
    @(184,28,'48'),
    @(185,28,'32'),
    @(186,28,'40'),
    @(187,28,'515'),
    @(188,28,'17'),
    @(189,28,'54200'),
    @(190,28,'37500'),
    @(191,28,'16900'),
    @(192,28,'9400'),
    @(193,28,'43300'),
    @(194,28,'1600'),
    @(195,28,'5600'),
    @(196,28,'49000'),
    @(206,28,'60.0'),
    @(210,28,'121.0')
)

foreach ($edit in $numericTextEdits) {
    $cell = $ws.Cells.Item($edit[0], $edit[1])
    $cell.NumberFormat = "@"
    $cell.Value = $edit[2]
}
